$wb = $excel.ActiveWorkbook

# Delete the "Datos Evaluador" worksheet entirely.
$excel.DisplayAlerts = $false
$wsDatos = $wb.Worksheets.Item("Datos Evaluador")
$wsDatos.Delete()
$excel.DisplayAlerts = $true

# Adjust selections to match the target state.
$wsF1 = $wb.Worksheets.Item("Formulario 1")
$wsF1.Range("C7:F7").Select()

$wsF3 = $wb.Worksheets.Item("Formulario 3")
$wsF4 = $wb.Worksheets.Item("Formulario 4")

$wsF4.Range("J13").Select()
$wsF4.Activate()
